$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.179.70'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.249.86'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.57%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.624'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '77.25'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.45%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.623'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0956'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.04'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.95%  '
$ws.Range("E13").Value = '  -3.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.584.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.86'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.858'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.254.70'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.081.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.82%  '
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("E20").Value = '  -2.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.02%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.33%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.25%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.41%  '
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("E34").Value = '  -4.71%  '
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("E36").Value = '  -1.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.42%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0303'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.29'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.89'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("E41").Value = '  -6.21%  '
$ws.Range("E42").Value = '  +12.23%  '
$ws.Range("E43").Value = '  -5.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.08%  '
$ws.Range("E46").Value = '  -3.19%  '
$ws.Range("E47").Value = '  -0.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.87%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.07%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.76%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.443'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +15.73%  '
